$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "64.398.71"
$ws.Range("E2").Value = "  +2.07%  "

# Row 3
Set-TextValue $ws "D3" "2.676.16"
$ws.Range("E3").Value = "  +3.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue $ws "D5" "596.42"
$ws.Range("E5").Value = "  +2.16%  "

# Row 6
Set-TextValue $ws "D6" "148.00"
$ws.Range("E6").Value = "  -0.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("E8").Value = "  -1.03%  "

# Row 9
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("E10").Value = "  -0.22%  "

# Row 12
Set-TextValue $ws "D12" "0.358"
$ws.Range("E12").Value = "  +1.03%  "

# Row 13
Set-TextValue $ws "D13" "27.96"
$ws.Range("E13").Value = "  +2.51%  "

# Row 14
Set-TextValue $ws "D14" "3.155.70"
$ws.Range("E14").Value = "  +3.00%  "

# Row 15
Set-TextValue $ws "D15" "64.313.79"
$ws.Range("E15").Value = "  +2.13%  "

# Row 16
$ws.Range("E16").Value = "  +0.42%  "

# Row 17
Set-TextValue $ws "D17" "2.713.48"
$ws.Range("E17").Value = "  +4.46%  "

# Row 18
$ws.Range("E18").Value = "  +0.65%  "

# Row 19
Set-TextValue $ws "D19" "346.57"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20
Set-TextValue $ws "D20" "4.41"
$ws.Range("E20").Value = "  +0.07%  "

# Row 21
Set-TextValue $ws "D21" "6.89"
$ws.Range("E21").Value = "  +1.48%  "

# Row 22
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
Set-TextValue $ws "D23" "68.79"
$ws.Range("E23").Value = "  +2.33%  "

# Row 24
$ws.Range("E24").Value = "  +11.02%  "

# Row 25
$ws.Range("E25").Value = "  +4.60%  "

# Row 26
$ws.Range("E26").Value = "  -1.11%  "

# Row 27
Set-TextValue $ws "D27" "8.55"
$ws.Range("E27").Value = "  +1.69%  "

# Row 28
Set-TextValue $ws "D28" "8.01"
$ws.Range("E28").Value = "  +1.15%  "

# Row 29
Set-TextValue $ws "D29" "1.00"
$ws.Range("E29").Value = "  +0.26%  "

# Row 30
Set-TextValue $ws "D30" "529.63"
$ws.Range("E30").Value = "  +13.66%  "

# Row 31
$ws.Range("E31").Value = "  +3.56%  "

# Row 32
Set-TextValue $ws "D32" "1.80"
$ws.Range("E32").Value = "  +11.90%  "

# Row 33
$subThree = [string][char]0x2083
Set-TextValue $ws "D33" ("0.0" + $subThree + "0829")
$ws.Range("E33").Value = "  +0.74%  "

# Row 34
Set-TextValue $ws "D34" "175.76"
$ws.Range("E34").Value = "  -0.61%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("E36").Value = "  +0.32%  "

# Row 37
$ws.Range("E37").Value = "  +0.57%  "

# Row 38
Set-TextValue $ws "D38" "4.69"
$ws.Range("E38").Value = "  +1.96%  "

# Row 39
Set-TextValue $ws "D39" "1.77"
$ws.Range("E39").Value = "  +3.63%  "

# Row 40
Set-TextValue $ws "D40" "172.68"
$ws.Range("E40").Value = "  +8.63%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
Set-TextValue $ws "D42" "40.76"
$ws.Range("E42").Value = "  +3.12%  "

# Row 43
$ws.Range("E43").Value = "  -0.29%  "

# Row 44
Set-TextValue $ws "D44" "21.83"
$ws.Range("E44").Value = "  +3.26%  "

# Row 45
$ws.Range("E45").Value = "  -0.49%  "

# Row 46
$ws.Range("E46").Value = "  +0.68%  "

# Row 47
$ws.Range("E47").Value = "  +1.59%  "

# Row 48
Set-TextValue $ws "D48" "0.0964"
$ws.Range("E48").Value = "  -1.03%  "

# Row 49
Set-TextValue $ws "D49" "18.87"
$ws.Range("E49").Value = "  +1.56%  "

# Row 50
$ws.Range("E50").Value = "  +2.82%  "

# Row 51
$ws.Range("E51").Value = "  -0.50%  "
